# Edits to distance descriptor
#
# 1) Resize three labeled callout text boxes ("Height bar", "Camera Holder",
#    "Scaffold Support"): widen/heighten them and make their text bold, 15pt.
# 2) Nudge the "Scaffold Support" leader-line connector's start point left to
#    follow the widened box.
# 3) Re-point the presentation's active theme palette ("simple-light-2") at
#    the alternate "Custom Theme" palette that previously only backed the
#    notes master (swapping which palette is the visible/active one).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function ConvertTo-RGBInt([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# --- "Height bar" callout (Shape 56): widen 1032300 -> 1130700 EMU; bold 15pt ---
$heightBar = $s.Shapes.Item("Shape 56")
$heightBar.Width = 89.03149796299212
$heightBar.TextFrame.TextRange.Font.Bold = $true
$heightBar.TextFrame.TextRange.Font.Size = 15

# --- "Camera Holder" callout (Shape 58): narrow 1032300 -> 888600 EMU; bold 15pt ---
$cameraHolder = $s.Shapes.Item("Shape 58")
$cameraHolder.Width = 69.96850393700787
$cameraHolder.TextFrame.TextRange.Font.Bold = $true
$cameraHolder.TextFrame.TextRange.Font.Size = 15

# --- "Scaffold Support" callout (Shape 59): 1032300x542400 -> 958500x642300 EMU; bold 15pt ---
$scaffoldSupport = $s.Shapes.Item("Shape 59")
$scaffoldSupport.Width = 75.47244094488188
$scaffoldSupport.Height = 50.5748043496063
$scaffoldSupport.TextFrame.TextRange.Font.Bold = $true
$scaffoldSupport.TextFrame.TextRange.Font.Size = 15

# --- Leader connector into "Scaffold Support" (Shape 60): off x 3146275 -> 3109375 EMU ---
$scaffoldConnector = $s.Shapes.Item("Shape 60")
$scaffoldConnector.Left = 244.83267976535433

# --- Swap the active theme palette over to the "Custom Theme" colors ---
$themeColors = $s.ThemeColorScheme
$themeColors.Item(1).RGB  = ConvertTo-RGBInt "000000"  # dk1
$themeColors.Item(2).RGB  = ConvertTo-RGBInt "FFFFFF"  # lt1
$themeColors.Item(3).RGB  = ConvertTo-RGBInt "158158"  # dk2
$themeColors.Item(4).RGB  = ConvertTo-RGBInt "F3F3F3"  # lt2
$themeColors.Item(5).RGB  = ConvertTo-RGBInt "058DC7"  # accent1
$themeColors.Item(6).RGB  = ConvertTo-RGBInt "50B432"  # accent2
$themeColors.Item(7).RGB  = ConvertTo-RGBInt "ED561B"  # accent3
$themeColors.Item(8).RGB  = ConvertTo-RGBInt "EDEF00"  # accent4
$themeColors.Item(9).RGB  = ConvertTo-RGBInt "24CBE5"  # accent5
$themeColors.Item(10).RGB = ConvertTo-RGBInt "64E572"  # accent6
$themeColors.Item(11).RGB = ConvertTo-RGBInt "2200CC"  # hlink
$themeColors.Item(12).RGB = ConvertTo-RGBInt "551A8B"  # folHlink
